$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.597.64"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.564.38"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'210.65"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "'0.487"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'24.89"
$ws.Range("E8").Value = "  +5.07%  "
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "'0.0584"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.789.63"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.555.41"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "28.639.84"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "'3.62"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "'61.28"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'230.94"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'150.54"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "'14.77"
$ws.Range("E26").Value = "  -1.02%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "'6.21"
$ws.Range("E30").Value = "  -4.61%  "
$ws.Range("D31").Value = "'1.06"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "1.389.44"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  -4.04%  "
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").Value = "'2.70"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "'0.0462"
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("D45").Value = "'63.72"
$ws.Range("E45").Value = "  +2.21%  "
$ws.Range("D46").Value = "'5.23"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("D47").Value = "1.701.26"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("E48").Value = "  -5.32%  "
$ws.Range("D49").Value = "'85.18"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'43.13"
$ws.Range("E50").Value = "  +6.11%  "
$ws.Range("D51").Value = "0.0₆0101"
$ws.Range("E51").Value = "  +2.00%  "
